$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (values that will not be misread as numbers)
$ws.Cells.Item(2, 4).Value = '28.848.74'
$ws.Cells.Item(3, 4).Value = '1.832.60'
$ws.Cells.Item(3, 5).Value = '  -1.73%  '
$ws.Cells.Item(4, 5).Value = '  -0.08%  '
$ws.Cells.Item(5, 5).Value = '  +0.65%  '
$ws.Cells.Item(6, 5).Value = '  -1.81%  '
$ws.Cells.Item(7, 5).Value = '  -0.06%  '
$ws.Cells.Item(8, 5).Value = '  -2.59%  '
$ws.Cells.Item(9, 5).Value = '  -2.64%  '
$ws.Cells.Item(10, 5).Value = '  -4.28%  '
$ws.Cells.Item(11, 5).Value = '  +0.29%  '
$ws.Cells.Item(12, 4).Value = '1.835.53'
$ws.Cells.Item(12, 5).Value = '  -1.67%  '
$ws.Cells.Item(13, 2).Value = 'Litecoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(13, 5).Value = '  +0.32%  '
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 5).Value = '  -0.99%  '
$ws.Cells.Item(15, 5).Value = '  -2.60%  '
$ws.Cells.Item(16, 5).Value = '  -1.14%  '
$ws.Cells.Item(17, 5).Value = '  -3.57%  '
$ws.Cells.Item(18, 4).Value = '28.857.97'
$ws.Cells.Item(18, 5).Value = '  -1.80%  '
$ws.Cells.Item(19, 5).Value = '  -2.51%  '
$ws.Cells.Item(20, 4).Value = '2.073.29'
$ws.Cells.Item(20, 5).Value = '  -2.30%  '
$ws.Cells.Item(21, 5).Value = '  -2.18%  '
$ws.Cells.Item(22, 5).Value = '  +0.09%  '
$ws.Cells.Item(23, 5).Value = '  -1.79%  '
$ws.Cells.Item(24, 5).Value = '  -0.03%  '
$ws.Cells.Item(25, 5).Value = '  -4.06%  '
$ws.Cells.Item(26, 5).Value = '  -0.96%  '
$ws.Cells.Item(27, 5).Value = '  -2.04%  '
$ws.Cells.Item(28, 5).Value = '  -2.43%  '
$ws.Cells.Item(29, 5).Value = '  -2.16%  '
$ws.Cells.Item(30, 5).Value = '  -1.64%  '
$ws.Cells.Item(31, 5).Value = '  -2.07%  '
$ws.Cells.Item(32, 5).Value = '  -0.20%  '
$ws.Cells.Item(33, 5).Value = '  -3.17%  '
$ws.Cells.Item(34, 5).Value = '  +2.64%  '
$ws.Cells.Item(35, 5).Value = '  -1.84%  '
$ws.Cells.Item(36, 5).Value = '  -3.35%  '
$ws.Cells.Item(37, 5).Value = '  -0.59%  '
$ws.Cells.Item(38, 5).Value = '  -0.74%  '
$ws.Cells.Item(39, 4).Value = '1.238.26'
$ws.Cells.Item(39, 5).Value = '  -2.71%  '
$ws.Cells.Item(40, 5).Value = '  -2.05%  '
$ws.Cells.Item(41, 5).Value = '  +6.35%  '
$ws.Cells.Item(42, 5).Value = '  -1.43%  '
$ws.Cells.Item(43, 5).Value = '  -0.12%  '
$ws.Cells.Item(44, 5).Value = '  +0.00%  '
$ws.Cells.Item(45, 5).Value = '  +0.57%  '
$ws.Cells.Item(46, 5).Value = '  -2.10%  '
$ws.Cells.Item(47, 2).Value = 'Mantle'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(47, 5).Value = '  -0.31%  '
$ws.Cells.Item(48, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(48, 4).Value = '1.975.46'
$ws.Cells.Item(48, 5).Value = '  -2.21%  '
$ws.Cells.Item(49, 5).Value = '  -8.99%  '
$ws.Cells.Item(50, 5).Value = '  -2.75%  '
$ws.Cells.Item(51, 5).Value = '  -1.39%  '

# Updates whose new value looks numeric -- force Text storage to match
# the original inlineStr/text cell type (avoid Excel auto-converting to Number).
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '244.51'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.6888'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.07690'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.3042'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '23.33'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.07810'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '92.63'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.093'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.6795'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.444'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.000008289'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '242.24'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '12.70'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.439'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '158.83'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '8.785'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '18.24'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.541'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.219'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.154'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.05081'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.7766'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.852'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.142'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.697'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.01850'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.694'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9535'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '108.19'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.937'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '9.628'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.00000000123'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.5158'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '63.86'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.742'
$cell.Style = 'Normal'
